# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The sheet's "K" column (column G, header "K" in G1) is recomputed from the
# refreshed strike data ("Strike#" source replaced). This writes the newly
# computed K values for every data row (rows 2-76), leaving all other
# columns (A date index, B date, C TB, D PC, E dS0, F dSF, H IP, I I0, J IF)
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K (column G) values for rows 2..76, in order.
$kValues = @(
    2,2,1,0,0,0,1,0,0,0,
    0,1,0,0,1,0,1,1,1,2,
    2,0,3,0,1,1,0,3,1,0,
    0,1,2,1,1,1,1,0,3,0,
    0,2,1,0,0,0,0,0,2,0,
    1,2,0,0,1,1,0,0,1,0,
    0,0,2,1,2,2,0,1,2,2,
    2,2,0,1,1
)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
